$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.072131
$ws.Range("H2").Value = 18.216393
$ws.Range("I2").Value = 0.003943999267036455
$ws.Range("J2").Value = 0.003943999267036454
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.2763116666666667
$ws.Range("N2").Value = 0.828935
$ws.Range("O2").Value = 0.02083107478128044
$ws.Range("P2").Value = 0.02083107478128044
$ws.Range("Q2").Value = 1.677800636828333
$ws.Range("R2").Value = 15.100205731455
$ws.Range("S2").Value = 0.00008215774366895162
$ws.Range("T2").Value = 0.00008215774366895161

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.072131
$ws.Range("H3").Value = 18.216393
$ws.Range("I3").Value = 0.003943999267036455
$ws.Range("J3").Value = 0.003943999267036454
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.180798333333333
$ws.Range("N3").Value = 6.542395
$ws.Range("O3").Value = 0.1644098988384798
$ws.Range("P3").Value = 0.1644098988384798
$ws.Range("Q3").Value = 13.24209316458166
$ws.Range("R3").Value = 119.178838481235
$ws.Range("S3").Value = 0.0006484325205125018
$ws.Range("T3").Value = 0.0006484325205125017

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.072131
$ws.Range("H4").Value = 18.216393
$ws.Range("I4").Value = 0.003943999267036455
$ws.Range("J4").Value = 0.003943999267036454
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 10.807288
$ws.Range("N4").Value = 32.421864
$ws.Range("O4").Value = 0.8147590263802398
$ws.Range("P4").Value = 0.8147590263802398
$ws.Range("Q4").Value = 65.623268490728
$ws.Range("R4").Value = 590.609416416552
$ws.Range("S4").Value = 0.003213409002855001
$ws.Range("T4").Value = 0.003213409002855

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1480.851806666667
$ws.Range("H5").Value = 4442.55542
$ws.Range("I5").Value = 0.9618498744646554
$ws.Range("J5").Value = 0.9618498744646552
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.2763116666666667
$ws.Range("N5").Value = 0.828935
$ws.Range("O5").Value = 0.02083107478128044
$ws.Range("P5").Value = 0.02083107478128044
$ws.Range("Q5").Value = 409.1766307864111
$ws.Range("R5").Value = 3682.589677077699
$ws.Range("S5").Value = 0.02003636666333844
$ws.Range("T5").Value = 0.02003636666333843

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1480.851806666667
$ws.Range("H6").Value = 4442.55542
$ws.Range("I6").Value = 0.9618498744646554
$ws.Range("J6").Value = 0.9618498744646552
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.180798333333333
$ws.Range("N6").Value = 6.542395
$ws.Range("O6").Value = 0.1644098988384798
$ws.Range("P6").Value = 0.1644098988384798
$ws.Range("Q6").Value = 3229.439151892322
$ws.Range("R6").Value = 29064.9523670309
$ws.Range("S6").Value = 0.1581376405585385
$ws.Range("T6").Value = 0.1581376405585384

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1480.851806666667
$ws.Range("H7").Value = 4442.55542
$ws.Range("I7").Value = 0.9618498744646554
$ws.Range("J7").Value = 0.9618498744646552
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 10.807288
$ws.Range("N7").Value = 32.421864
$ws.Range("O7").Value = 0.8147590263802398
$ws.Range("P7").Value = 0.8147590263802398
$ws.Range("Q7").Value = 16003.99195996698
$ws.Range("R7").Value = 144035.9276397029
$ws.Range("S7").Value = 0.7836758672427785
$ws.Range("T7").Value = 0.7836758672427784

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 52.663316
$ws.Range("H8").Value = 157.989948
$ws.Range("I8").Value = 0.03420612626830831
$ws.Range("J8").Value = 0.0342061262683083
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.2763116666666667
$ws.Range("N8").Value = 0.828935
$ws.Range("O8").Value = 0.02083107478128044
$ws.Range("P8").Value = 0.02083107478128044
$ws.Range("Q8").Value = 14.55148861615333
$ws.Range("R8").Value = 130.96339754538
$ws.Range("S8").Value = 0.0007125503742730515
$ws.Range("T8").Value = 0.0007125503742730514

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 52.663316
$ws.Range("H9").Value = 157.989948
$ws.Range("I9").Value = 0.03420612626830831
$ws.Range("J9").Value = 0.0342061262683083
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.180798333333333
$ws.Range("N9").Value = 6.542395
$ws.Range("O9").Value = 0.1644098988384798
$ws.Range("P9").Value = 0.1644098988384798
$ws.Range("Q9").Value = 114.8480717606067
$ws.Range("R9").Value = 1033.63264584546
$ws.Range("S9").Value = 0.005623825759428834
$ws.Range("T9").Value = 0.005623825759428833

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 52.663316
$ws.Range("H10").Value = 157.989948
$ws.Range("I10").Value = 0.03420612626830831
$ws.Range("J10").Value = 0.0342061262683083
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 10.807288
$ws.Range("N10").Value = 32.421864
$ws.Range("O10").Value = 0.8147590263802398
$ws.Range("P10").Value = 0.8147590263802398
$ws.Range("Q10").Value = 569.147623047008
$ws.Range("R10").Value = 5122.328607423072
$ws.Range("S10").Value = 0.02786975013460642
$ws.Range("T10").Value = 0.02786975013460642
